$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E data rows 2-9: numeric 83 -> text "November 10 - November 16"
$ws.Range("E2:E9").Value = "November 10 - November 16"

# Column E header: "Occurrence " -> "Dates Used"
$ws.Range("E1").Value = "Dates Used"

# Reset the active selection to E1 (matches the sheet view state after the edit)
$ws.Range("E1").Select()
